$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "g vs. hardcoded C 1 checkers"

# Update win rate (B) and avg. game length (C) values per row
$ws.Range("B2").Value = 34
$ws.Range("C2").Value = 12.821
$ws.Range("B3").Value = 82
$ws.Range("C3").Value = 4.965
$ws.Range("C4").Value = 22.352
$ws.Range("B5").Value = 67.1
$ws.Range("C5").Value = 5.227
$ws.Range("B6").Value = 84.7
$ws.Range("C6").Value = 4.205
$ws.Range("C7").Value = 7.013
$ws.Range("B8").Value = 85.3
$ws.Range("C8").Value = 9.53
$ws.Range("B9").Value = 48.6
$ws.Range("C9").Value = 4.011
$ws.Range("C12").Value = 8.801
$ws.Range("B13").Value = 31.9
$ws.Range("C13").Value = 22.788
$ws.Range("B14").Value = 52.1
$ws.Range("C14").Value = 26.748
$ws.Range("B15").Value = 16.7
$ws.Range("C15").Value = 22.167
$ws.Range("B16").Value = 83.7
$ws.Range("C16").Value = 6.582
$ws.Range("B17").Value = 35.2
$ws.Range("C17").Value = 9.044
$ws.Range("C18").Value = 5.169
$ws.Range("C19").Value = 14.958
$ws.Range("C21").Value = 13.446
$ws.Range("C22").Value = 22.507
$ws.Range("C26").Value = 15.823
$ws.Range("C27").Value = 8.96
$ws.Range("B28").Value = 85.3
$ws.Range("C28").Value = 5.58
$ws.Range("B29").Value = 66.6
$ws.Range("C29").Value = 4.041
$ws.Range("C30").Value = 16.479
$ws.Range("B31").Value = 19.6
$ws.Range("C31").Value = 27.12
$ws.Range("B32").Value = 47.3
$ws.Range("C32").Value = 8.299
$ws.Range("C33").Value = 5.021
$ws.Range("C34").Value = 4.642
$ws.Range("B35").Value = 83.9
$ws.Range("C35").Value = 16.872
$ws.Range("B36").Value = 49.7
$ws.Range("C36").Value = 3.627
$ws.Range("C37").Value = 17.92
$ws.Range("B38").Value = 80.8
$ws.Range("C38").Value = 7.204
$ws.Range("B40").Value = 49.6
$ws.Range("C40").Value = 8.792
$ws.Range("C41").Value = 5
$ws.Range("B42").Value = 49.3
$ws.Range("C42").Value = 3.862
$ws.Range("B45").Value = 17.4
$ws.Range("C45").Value = 13.545
$ws.Range("B46").Value = 51.3
$ws.Range("C46").Value = 18.642
$ws.Range("B47").Value = 15.8
$ws.Range("C47").Value = 4.466
